# Update the workout/activity log on Sheet1:
#   - add an "Intensity Rating (1-10)" column (E)
#   - split the old "Gym" entries into separate "Strength Training - Chest"
#     and "Strength Training - Legs" rows
#   - rename "Duration" -> "Duration (min)" and "Distance" -> "Distance (m)"
#   - fill in the new Duration/Distance/Intensity figures
#   - widen the columns to fit the longer header/activity text
#
# NOTE: values are written in the same order the shared-string table in the
# target workbook implies (new intensity header, then the two strength
# training labels, then the renamed duration/distance headers) so the
# resulting xl/sharedStrings.xml ordering lines up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Intensity Rating (1-10)" header (column E) ---
$ws.Range("E1").Value = "Intensity Rating (1-10)"

# --- Split "Gym" into two specific strength-training activities ---
$ws.Range("B4").Value = "Strength Training - Chest"
$ws.Range("B5").Value = "Strength Training - Legs"

# --- Rename existing headers to include units ---
$ws.Range("C1").Value = "Duration (min)"
$ws.Range("D1").Value = "Distance (m)"

# --- Row 2: Swimming ---
$ws.Range("D2").Value = 550
$ws.Range("E2").Value = 8

# --- Row 4: Strength Training - Chest ---
$ws.Range("C4").Value = 90
$ws.Range("E4").Value = 8

# --- Row 5: Strength Training - Legs ---
$ws.Range("C5").Value = 120
$ws.Range("E5").Value = 7

# --- Row 6: Swimming ---
$ws.Range("D6").Value = 650
$ws.Range("E6").Value = 7

# --- Row 7: Running ---
$ws.Range("D7").Value = 6470
$ws.Range("E7").Value = 8

# --- Column widths so the new/renamed headers are readable ---
# (ColumnWidth is quantized to Excel's pixel grid by the host, so these
# inputs land on the closest attainable width to the authored values.)
$ws.Columns.Item(2).ColumnWidth = 23.3
$ws.Columns.Item(3).ColumnWidth = 18.33
$ws.Columns.Item(4).ColumnWidth = 10.33
$ws.Columns.Item(5).ColumnWidth = 19.67

# --- Selection, matching the author's last cursor position ---
$ws.Range("I17:I18").Select() | Out-Null
